$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 3642.6
$ws.Range("I38").Value = 19
$ws.Range("J38").Value = 7783.857
$ws.Range("K38").Value = 57
$ws.Range("L38").Value = 23351.571
$ws.Range("M38").Value = 315
$ws.Range("N38").Value = -24095.571
# Row 61
$ws.Range("H61").Value = 233.33333
$ws.Range("I61").Value = 233.33333
$ws.Range("K61").Value = 699.99999
$ws.Range("M61").Value = -527.99999
# Row 100
$ws.Range("H100").Value = 3444.3572
$ws.Range("I100").Value = 3479.111
$ws.Range("K100").Value = 3479.111
$ws.Range("M100").Value = -2938.111
# Row 106
$ws.Range("H106").Value = 7111.1113
$ws.Range("I106").Value = 8083.3335
$ws.Range("J106").Value = 5166.6665
$ws.Range("K106").Value = 8083.3335
$ws.Range("L106").Value = 5166.6665
$ws.Range("M106").Value = -7452.3335
$ws.Range("N106").Value = -6428.6665
# Row 112
$ws.Range("H112").Value = 1526.6471
$ws.Range("I112").Value = 1489
$ws.Range("J112").Value = 1529
$ws.Range("K112").Value = 4467
$ws.Range("L112").Value = 4587
$ws.Range("M112").Value = -3359
$ws.Range("N112").Value = -6803
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6241.516
$ws.Range("I32").Value = 4183.375
$ws.Range("K32").Value = 4183.375
$ws.Range("M32").Value = -3896.375
# Row 45
$ws.Range("H45").Value = 1877.96
$ws.Range("I45").Value = 1719.1052
$ws.Range("J45").Value = 2381
$ws.Range("K45").Value = 1719.1052
$ws.Range("L45").Value = 2381
$ws.Range("M45").Value = -1342.1052
$ws.Range("N45").Value = -3135
# Row 61
$ws.Range("H61").Value = 13764.821
$ws.Range("I61").Value = 2375.6191
$ws.Range("J61").Value = 47932.43
$ws.Range("K61").Value = 2375.6191
$ws.Range("L61").Value = 47932.43
$ws.Range("M61").Value = -2163.6191
$ws.Range("N61").Value = -48356.43
# Row 74
$ws.Range("H74").Value = 76679.21000000001
$ws.Range("I74").Value = 51310.5
$ws.Range("J74").Value = 152785.33
$ws.Range("K74").Value = 51310.5
$ws.Range("L74").Value = 152785.33
$ws.Range("M74").Value = -50436.5
$ws.Range("N74").Value = -154533.33
# Row 77
$ws.Range("H77").Value = 76679.21000000001
$ws.Range("I77").Value = 51310.5
$ws.Range("J77").Value = 152785.33
$ws.Range("K77").Value = 256552.5
$ws.Range("L77").Value = 763926.6499999999
$ws.Range("M77").Value = -252184.5
$ws.Range("N77").Value = -772662.6499999999
# Row 122
$ws.Range("H122").Value = 132126.58
$ws.Range("I122").Value = 2295.6667
$ws.Range("J122").Value = 229499.75
$ws.Range("K122").Value = 6887.000100000001
$ws.Range("L122").Value = 688499.25
$ws.Range("M122").Value = -4437.000100000001
$ws.Range("N122").Value = -693399.25
# Row 132
$ws.Range("H132").Value = 8550.611999999999
$ws.Range("I132").Value = 8932.139999999999
$ws.Range("J132").Value = 5816.3335
$ws.Range("K132").Value = 26796.42
$ws.Range("L132").Value = 17449.0005
$ws.Range("M132").Value = -24266.42
$ws.Range("N132").Value = -22509.0005
# Row 136
$ws.Range("H136").Value = 13764.821
$ws.Range("I136").Value = 2375.6191
$ws.Range("J136").Value = 47932.43
$ws.Range("K136").Value = 7126.8573
$ws.Range("L136").Value = 143797.29
$ws.Range("M136").Value = -4576.8573
$ws.Range("N136").Value = -148897.29
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 29615.355
$ws.Range("I94").Value = 501.43478
$ws.Range("K94").Value = 501.43478
$ws.Range("M94").Value = -50.43477999999999
# Row 99
$ws.Range("H99").Value = 2167.125
$ws.Range("I99").Value = 1393.1177
$ws.Range("J99").Value = 4046.8572
$ws.Range("K99").Value = 1393.1177
$ws.Range("L99").Value = 4046.8572
$ws.Range("M99").Value = 104.8823
$ws.Range("N99").Value = -7042.8572
# Row 134
$ws.Range("H134").Value = 3300.25
$ws.Range("I134").Value = 2277.1667
$ws.Range("K134").Value = 6831.500100000001
$ws.Range("M134").Value = -4296.500100000001
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2165.973
$ws.Range("I31").Value = 1282.7333
$ws.Range("J31").Value = 2768.182
$ws.Range("K31").Value = 1282.7333
$ws.Range("L31").Value = 2768.182
$ws.Range("M31").Value = -987.7333000000001
$ws.Range("N31").Value = -3358.182
# Row 34
$ws.Range("H34").Value = 2165.973
$ws.Range("I34").Value = 1282.7333
$ws.Range("J34").Value = 2768.182
$ws.Range("K34").Value = 1282.7333
$ws.Range("L34").Value = 2768.182
$ws.Range("M34").Value = -1080.7333
$ws.Range("N34").Value = -3172.182
# Row 58
$ws.Range("H58").Value = 4932.263
$ws.Range("I58").Value = 4719.933
$ws.Range("K58").Value = 4719.933
$ws.Range("M58").Value = -4516.933
# Row 132
$ws.Range("H132").Value = 3860.8206
$ws.Range("I132").Value = 1175
$ws.Range("J132").Value = 16138.857
$ws.Range("K132").Value = 3525
$ws.Range("L132").Value = 48416.571
$ws.Range("M132").Value = -995
$ws.Range("N132").Value = -53476.571
# Row 134
$ws.Range("H134").Value = 2366.4915
$ws.Range("I134").Value = 2188.5615
$ws.Range("J134").Value = 7437.5
$ws.Range("K134").Value = 6565.684499999999
$ws.Range("L134").Value = 22312.5
$ws.Range("M134").Value = -4030.684499999999
$ws.Range("N134").Value = -27382.5
# Row 135
$ws.Range("H135").Value = 67999.8
$ws.Range("J135").Value = 67999.8
$ws.Range("L135").Value = 67999.8
$ws.Range("N135").Value = -78139.8
# Row 136
$ws.Range("H136").Value = 4932.263
$ws.Range("I136").Value = 4719.933
$ws.Range("K136").Value = 14159.799
$ws.Range("M136").Value = -11609.799
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 186026.36
$ws.Range("I4").Value = 202503.4
$ws.Range("J4").Value = 10957.625
$ws.Range("K4").Value = 607510.2
$ws.Range("L4").Value = 32872.875
$ws.Range("M4").Value = -607398.2
$ws.Range("N4").Value = -33096.875
# Row 9
$ws.Range("H9").Value = 19275.125
$ws.Range("I9").Value = 30240.2
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 90720.60000000001
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = -90496.60000000001
$ws.Range("N9").Value = -3448
# Row 131
$ws.Range("H131").Value = 3137.7144
$ws.Range("I131").Value = 1960.5
$ws.Range("J131").Value = 4207.909
$ws.Range("K131").Value = 5881.5
$ws.Range("L131").Value = 12623.727
$ws.Range("M131").Value = -841.5
$ws.Range("N131").Value = -22703.727
$ws = $wb.Worksheets.Item("GSM")
# Row 33
$ws.Range("H33").Value = 43125
$ws.Range("I33").Value = 39500
$ws.Range("K33").Value = 39500
$ws.Range("M33").Value = -39248
# Row 97
$ws.Range("H97").Value = 16379.5
$ws.Range("I97").Value = 23422.293
$ws.Range("J97").Value = 1181.8948
$ws.Range("K97").Value = 23422.293
$ws.Range("L97").Value = 1181.8948
$ws.Range("M97").Value = -22926.293
$ws.Range("N97").Value = -2173.8948
# Row 102
$ws.Range("H102").Value = 21805.068
$ws.Range("I102").Value = 24011.305
$ws.Range("J102").Value = 13998.385
$ws.Range("K102").Value = 24011.305
$ws.Range("L102").Value = 13998.385
$ws.Range("M102").Value = -22389.305
$ws.Range("N102").Value = -17242.385
# Row 113
$ws.Range("H113").Value = 4129.0625
$ws.Range("I113").Value = 4460.3687
$ws.Range("J113").Value = 3644.8462
$ws.Range("K113").Value = 4460.3687
$ws.Range("L113").Value = 3644.8462
$ws.Range("M113").Value = -2290.3687
$ws.Range("N113").Value = -7984.8462
# Row 122
$ws.Range("H122").Value = 45922.9
$ws.Range("I122").Value = 55296.938
$ws.Range("K122").Value = 165890.814
$ws.Range("M122").Value = -163440.814
# Row 132
$ws.Range("H132").Value = 3247.4788
$ws.Range("I132").Value = 3551.698
$ws.Range("J132").Value = 2351.7222
$ws.Range("K132").Value = 10655.094
$ws.Range("L132").Value = 7055.1666
$ws.Range("M132").Value = -8125.093999999999
$ws.Range("N132").Value = -12115.1666
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 59524.938
$ws.Range("I136").Value = 3385.6428
$ws.Range("J136").Value = 452500
$ws.Range("K136").Value = 10156.9284
$ws.Range("L136").Value = 1357500
$ws.Range("M136").Value = -7606.928400000001
$ws.Range("N136").Value = -1362600
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 5214.815
$ws.Range("I122").Value = 5101.619
$ws.Range("J122").Value = 5611
$ws.Range("K122").Value = 15304.857
$ws.Range("L122").Value = 16833
$ws.Range("M122").Value = -12854.857
$ws.Range("N122").Value = -21733
# Row 132
$ws.Range("H132").Value = 3043.6965
$ws.Range("I132").Value = 2466
$ws.Range("J132").Value = 7087.5713
$ws.Range("K132").Value = 7398
$ws.Range("L132").Value = 21262.7139
$ws.Range("M132").Value = -4868
$ws.Range("N132").Value = -26322.7139
# Row 136
$ws.Range("H136").Value = 5284.6924
$ws.Range("I136").Value = 5653.815
$ws.Range("J136").Value = 4454.1665
$ws.Range("K136").Value = 16961.445
$ws.Range("L136").Value = 13362.4995
$ws.Range("M136").Value = -14411.445
$ws.Range("N136").Value = -18462.4995
